$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 1; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = "MS"
}

$ws.Range("I12").Select() | Out-Null
